# "excel, lvl 3 stage clear and reactivate"
#
# Updates the AFTER ALPHA task table on the "Full 1" sheet:
#  - row 59 (0.6 / -Fix all current bugs): Real time now also includes Martí's hour
#  - row 76 (0.8 / -//1/2 enemies (20-27), i.e. "Colliders level 3" area): Martí
#    joins Jorge/Gerard on the task, and his time is added to the Real time total
#  - row 81 (-1/2 Boss (first phase).): task is reactivated with an estimated
#    and real time for Martí
#  - row 92 (Martí's summary row): Est. Total Time / Real time totals refreshed
#    to reflect the above additions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K59").Value = "J: 2h G: 1h  M: 1h"

$ws.Range("I76").Value = "Jorge/Gerard/Martí"
$ws.Range("K76").Value = "G: 4h M: 1h 30 min"

$ws.Range("J81").Value = "5h"
$ws.Range("K81").Value = "7h 50min"

$ws.Range("J92").Value = "23h"
$ws.Range("K92").Value = "36h 20 min"
